# Lab 6 zip-file commit: insert a new "ADD EXAMPLE" slide right before the
# existing "Summary of types of request-response Web APIs" slide (slide 22),
# pushing it and the following slides down by one position.

$p = $ppt.ActivePresentation

# The "Summary of types..." slide is the 22nd slide in the deck today and
# uses the standard "Title and Content" layout - reuse that same layout for
# the new slide so it gets a title placeholder + a body/content placeholder.
$summarySlide = $p.Slides.Item(22)
$layout = $summarySlide.CustomLayout

# Insert the new slide at index 22 - this bumps the "Summary..." slide (and
# everything after it) down by one, exactly like the diff's reordered
# <p:sldId> list.
$newSlide = $p.Slides.AddSlide(22, $layout)

# Title placeholder.
$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "ADD EXAMPLE"

# Body/content placeholder.
$newSlide.Shapes.Item(2).TextFrame.TextRange.Text = "ADD EXAMPLE"
